# Restore project back to default: update the three sample e-mail addresses
# and move the active-cell selection from F7 to D7 (matches the
# "restored back to default" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "janfaizi1@gmail.com"
$ws.Range("C3").Value = "alijan2@tek.com"
$ws.Range("C4").Value = "anisa2@gmail.com"

$ws.Range("D7").Select()
